$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their text (string) representation instead of being
# auto-converted to numbers/percentages by Excel when the new value is set.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.20%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "47.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.32%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.075"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.06%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07704"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-4.75%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.72%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.291"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "18.08%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.559"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-7.49%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1232"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.79%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1923"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.68%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09228"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.47%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04523"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "6.54%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1048"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.30%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001285"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.44%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04207"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.68%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005909"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.02%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.339"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.80%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.401"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.39%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.02%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.164"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.75%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1394"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.69%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3036"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.54%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.36%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004105"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.62%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001361"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.12%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003571"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-95.18%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02563"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-5.19%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05666"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.76%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "84.89%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007960"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.38%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1418"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.47%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008401"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "14.46%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007770"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-9.50%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3384"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.73%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006832"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.30%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "1.15%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1.93%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004033"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "1.15%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002117"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "1.15%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002017"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "1.15%"
